$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 37794
$ws1.Range("F5").Value = 788
$ws1.Range("F9").Value = 862
$ws1.Range("F10").Value = 99
$ws1.Range("F11").Value = 734
$ws1.Range("F12").Value = 567
$ws1.Range("F13").Value = 67
$ws1.Range("F14").Value = 38
$ws1.Range("F15").Value = 30
$ws1.Range("F16").Value = 665
$ws1.Range("F17").Value = 185
$ws1.Range("F18").Value = 477
$ws1.Range("F20").Value = 1178
$ws1.Range("F21").Value = 95
$ws1.Range("F22").Value = 850
$ws1.Range("F23").Value = 2559
$ws1.Range("F24").Value = 1046
$ws1.Range("F27").Value = 1169
$ws1.Range("F29").Value = 805

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 422

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 649

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 649
$ws4.Range("F3").Value = 37794
$ws4.Range("F6").Value = 788
$ws4.Range("F11").Value = 422
$ws4.Range("F15").Value = 862
$ws4.Range("F16").Value = 99
$ws4.Range("F17").Value = 734
$ws4.Range("F18").Value = 567
$ws4.Range("F19").Value = 67
$ws4.Range("F21").Value = 38
$ws4.Range("F25").Value = 30
$ws4.Range("F27").Value = 665
$ws4.Range("F28").Value = 185
$ws4.Range("F29").Value = 477
$ws4.Range("F31").Value = 1178
$ws4.Range("F32").Value = 95
$ws4.Range("F33").Value = 850
$ws4.Range("F34").Value = 2559
$ws4.Range("F35").Value = 1046
$ws4.Range("F38").Value = 1169
$ws4.Range("F41").Value = 805
